$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3340876
$ws.Range("J17").Value = 3537104
$ws.Range("L17").Value = 10611312
$ws.Range("N17").Value = -10611648
$ws.Range("H129").Value = 1637.75
$ws.Range("I129").Value = 1010.3
$ws.Range("K129").Value = 3030.9
$ws.Range("M129").Value = 1969.1
$ws.Range("H132").Value = 8069.452
$ws.Range("I132").Value = 5959.325
$ws.Range("K132").Value = 17877.975
$ws.Range("M132").Value = -15347.975
$ws.Range("H135").Value = 2235.5
$ws.Range("I135").Value = 2450.7
$ws.Range("J135").Value = 1159.5
$ws.Range("K135").Value = 22056.3
$ws.Range("L135").Value = 10435.5
$ws.Range("M135").Value = -19521.3
$ws.Range("N135").Value = -15505.5
$ws.Range("H137").Value = 6085.0386
$ws.Range("I137").Value = 2421.4814
$ws.Range("J137").Value = 10041.68
$ws.Range("K137").Value = 7264.4442
$ws.Range("L137").Value = 30125.04
$ws.Range("M137").Value = -4714.4442
$ws.Range("N137").Value = -35225.04
$ws.Range("H138").Value = 2151.97
$ws.Range("I138").Value = 1163.5
$ws.Range("K138").Value = 3490.5
$ws.Range("M138").Value = 1649.5
$ws.Range("H141").Value = 3073.8293
$ws.Range("I141").Value = 2523
$ws.Range("J141").Value = 7039.8
$ws.Range("K141").Value = 7569
$ws.Range("L141").Value = 21119.4
$ws.Range("M141").Value = -2389
$ws.Range("N141").Value = -31479.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6781.154
$ws.Range("I32").Value = 5719.968
$ws.Range("K32").Value = 5719.968
$ws.Range("M32").Value = -5432.968
$ws.Range("H33").Value = 6378.625
$ws.Range("I33").Value = 3000
$ws.Range("K33").Value = 3000
$ws.Range("M33").Value = -2671
$ws.Range("H45").Value = 2902.451
$ws.Range("I45").Value = 2748.0488
$ws.Range("K45").Value = 2748.0488
$ws.Range("M45").Value = -2371.0488
$ws.Range("H61").Value = 168571.5
$ws.Range("I61").Value = 3992.25
$ws.Range("J61").Value = 234403.2
$ws.Range("K61").Value = 3992.25
$ws.Range("L61").Value = 234403.2
$ws.Range("M61").Value = -3780.25
$ws.Range("N61").Value = -234827.2
$ws.Range("H74").Value = 10690.6045
$ws.Range("I74").Value = 1422.8684
$ws.Range("K74").Value = 1422.8684
$ws.Range("M74").Value = -548.8684000000001
$ws.Range("H77").Value = 10690.6045
$ws.Range("I77").Value = 1422.8684
$ws.Range("K77").Value = 7114.342000000001
$ws.Range("M77").Value = -2746.342000000001
$ws.Range("H119").Value = 30000
$ws.Range("J119").Value = 30000
$ws.Range("L119").Value = 30000
$ws.Range("N119").Value = -39676
$ws.Range("H122").Value = 3751
$ws.Range("I122").Value = 1671.4706
$ws.Range("K122").Value = 5014.4118
$ws.Range("M122").Value = -2564.4118
$ws.Range("H136").Value = 168571.5
$ws.Range("I136").Value = 3992.25
$ws.Range("J136").Value = 234403.2
$ws.Range("K136").Value = 11976.75
$ws.Range("L136").Value = 703209.6000000001
$ws.Range("M136").Value = -9426.75
$ws.Range("N136").Value = -708309.6000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 20000
$ws.Range("J18").Value = 20000
$ws.Range("L18").Value = 20000
$ws.Range("N18").Value = -21058
$ws.Range("H22").Value = 1636.92
$ws.Range("I22").Value = 1887.4286
$ws.Range("J22").Value = 321.75
$ws.Range("K22").Value = 1887.4286
$ws.Range("L22").Value = 321.75
$ws.Range("M22").Value = -1714.4286
$ws.Range("N22").Value = -667.75
$ws.Range("H105").Value = 1294.1111
$ws.Range("I105").Value = 1294.1111
$ws.Range("K105").Value = 1294.1111
$ws.Range("M105").Value = 452.8888999999999
$ws.Range("H107").Value = 1601.76
$ws.Range("I107").Value = 1543.5
$ws.Range("K107").Value = 1543.5
$ws.Range("M107").Value = 376.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 957.14813
$ws.Range("I22").Value = 508.73685
$ws.Range("K22").Value = 508.73685
$ws.Range("M22").Value = -158.73685
$ws.Range("H31").Value = 22015.967
$ws.Range("I31").Value = 20656.25
$ws.Range("J31").Value = 23569.928
$ws.Range("K31").Value = 20656.25
$ws.Range("L31").Value = 23569.928
$ws.Range("M31").Value = -20361.25
$ws.Range("N31").Value = -24159.928
$ws.Range("H34").Value = 22015.967
$ws.Range("I34").Value = 20656.25
$ws.Range("J34").Value = 23569.928
$ws.Range("K34").Value = 20656.25
$ws.Range("L34").Value = 23569.928
$ws.Range("M34").Value = -20454.25
$ws.Range("N34").Value = -23973.928
$ws.Range("H39").Value = 8791.5
$ws.Range("I39").Value = 6500
$ws.Range("J39").Value = 9937.25
$ws.Range("K39").Value = 6500
$ws.Range("L39").Value = 9937.25
$ws.Range("M39").Value = -6109
$ws.Range("N39").Value = -10719.25
$ws.Range("H49").Value = 8791.5
$ws.Range("I49").Value = 6500
$ws.Range("J49").Value = 9937.25
$ws.Range("K49").Value = 6500
$ws.Range("L49").Value = 9937.25
$ws.Range("M49").Value = -6318
$ws.Range("N49").Value = -10301.25
$ws.Range("H62").Value = 4005
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H65").Value = 4005
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H105").Value = 8379.214
$ws.Range("I105").Value = 10931
$ws.Range("J105").Value = 1999.75
$ws.Range("K105").Value = 10931
$ws.Range("L105").Value = 1999.75
$ws.Range("M105").Value = -9184
$ws.Range("N105").Value = -5493.75
$ws.Range("H107").Value = 658.2
$ws.Range("I107").Value = 517.125
$ws.Range("J107").Value = 819.4286
$ws.Range("K107").Value = 517.125
$ws.Range("L107").Value = 819.4286
$ws.Range("M107").Value = 1402.875
$ws.Range("N107").Value = -4659.4286
$ws.Range("H132").Value = 8530.105
$ws.Range("I132").Value = 2141.818
$ws.Range("J132").Value = 17314
$ws.Range("K132").Value = 6425.454000000001
$ws.Range("L132").Value = 51942
$ws.Range("M132").Value = -3895.454000000001
$ws.Range("N132").Value = -57002
$ws.Range("H134").Value = 4852.2793
$ws.Range("I134").Value = 1807.0625
$ws.Range("K134").Value = 5421.1875
$ws.Range("M134").Value = -2886.1875
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1325
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 22317.92
$ws.Range("I70").Value = 21664.953
$ws.Range("K70").Value = 21664.953
$ws.Range("M70").Value = -21394.953
$ws.Range("H73").Value = 22317.92
$ws.Range("I73").Value = 21664.953
$ws.Range("K73").Value = 21664.953
$ws.Range("M73").Value = -20728.953
$ws.Range("H122").Value = 7714.7144
$ws.Range("I122").Value = 8500.5
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 25501.5
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -23051.5
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 9648.813
$ws.Range("I132").Value = 7330.3887
$ws.Range("J132").Value = 21572.143
$ws.Range("K132").Value = 21991.1661
$ws.Range("L132").Value = 64716.429
$ws.Range("M132").Value = -19461.1661
$ws.Range("N132").Value = -69776.429

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2141.8572
$ws.Range("I16").Value = 2257.7693
$ws.Range("J16").Value = 635
$ws.Range("K16").Value = 2257.7693
$ws.Range("L16").Value = 635
$ws.Range("M16").Value = -2087.7693
$ws.Range("N16").Value = -975
$ws.Range("H30").Value = 3393
$ws.Range("I30").Value = 2198.5
$ws.Range("J30").Value = 4189.3335
$ws.Range("K30").Value = 2198.5
$ws.Range("L30").Value = 4189.3335
$ws.Range("M30").Value = -2090.5
$ws.Range("N30").Value = -4405.3335
$ws.Range("H93").Value = 6479.077
$ws.Range("I93").Value = 4881.222
$ws.Range("J93").Value = 10074.25
$ws.Range("K93").Value = 4881.222
$ws.Range("L93").Value = 10074.25
$ws.Range("M93").Value = -3633.222
$ws.Range("N93").Value = -12570.25
$ws.Range("H106").Value = 19966.666
$ws.Range("J106").Value = 19966.666
$ws.Range("L106").Value = 19966.666
$ws.Range("N106").Value = -22490.666
$ws.Range("H136").Value = 65301.17
$ws.Range("I136").Value = 107188.48
$ws.Range("J136").Value = 15560
$ws.Range("K136").Value = 321565.44
$ws.Range("L136").Value = 46680
$ws.Range("M136").Value = -319015.44
$ws.Range("N136").Value = -51780
$ws.Range("H138").Value = 49999.5
$ws.Range("J138").Value = 49999.5
$ws.Range("L138").Value = 49999.5
$ws.Range("N138").Value = -60279.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 16713.125
$ws.Range("J74").Value = 16713.125
$ws.Range("L74").Value = 16713.125
$ws.Range("N74").Value = -18585.125
$ws.Range("H77").Value = 16713.125
$ws.Range("J77").Value = 16713.125
$ws.Range("L77").Value = 50139.375
$ws.Range("N77").Value = -59499.375
$ws.Range("H113").Value = 517.2
$ws.Range("I113").Value = 517.2
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1551.6
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 618.3999999999999
$ws.Range("H136").Value = 9529.950999999999
$ws.Range("I136").Value = 1074.2609
$ws.Range("K136").Value = 3222.7827
$ws.Range("M136").Value = -672.7826999999997
$ws.Range("N113").ClearContents()

